$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 4 / B4 ---------------------------------------------------------
# "1. User account`n2. Contract account" -> "1. Account`n2. Program"
# Font color also changes (red -> default/automatic). Copy formatting from
# B6 (already uses the target "no explicit color" style) before writing the
# new text so the engine reuses the existing cellXf instead of minting a
# brand-new one for an explicit black color.
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Value = "1. Account" + [char]10 + "2. Program"

# --- Row 5 / B5 ----------------------------------------------------------
# "nodeos" -> "test-ledger/" (same red -> default font-color change)
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Value = "test-ledger/"

# --- Row 11 / B11 ---------------------------------------------------------
# "uint32_t" -> "u32" (same red -> default font-color change)
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "u32"

# --- Row 17 / B17 ----------------------------------------------------------
# "EOSIO token" -> "Token program" (style/color unchanged - stays red)
$ws.Range("B17").Value = "Token program"

# Clear clipboard marching ants / selection artifacts from the Copy() calls
$excel.CutCopyMode = 0

# --- Selection / scroll position ------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B18").Select() | Out-Null
